$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  9"
$ws.Range("C9").Value = "Report Covering the Week  2/24/2025  Through  3/2/2025"

# --- Crime statistics data updates ---
# Row 14
$ws.Range("M14").Value = -50
$ws.Range("N14").Value = -80
# Row 15
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -20
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 80
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = -52.631578947368
# Row 16
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = -45.454545454545
$ws.Range("G16").Value = 27
$ws.Range("H16").Value = -18.518518518518
$ws.Range("I16").Value = 50
$ws.Range("J16").Value = 49
$ws.Range("K16").Value = 2.04081632653
$ws.Range("L16").Value = -24.242424242424
$ws.Range("M16").Value = -7.407407407407
$ws.Range("N16").Value = -80.237154150197
# Row 17
$ws.Range("C17").Value = 19
$ws.Range("D17").Value = 36
$ws.Range("E17").Value = -47.222222222222
$ws.Range("F17").Value = 66
$ws.Range("G17").Value = 96
$ws.Range("H17").Value = -31.25
$ws.Range("I17").Value = 133
$ws.Range("J17").Value = 160
$ws.Range("K17").Value = -16.875
$ws.Range("L17").Value = -5
$ws.Range("M17").Value = 66.25
$ws.Range("N17").Value = -21.764705882352
# Row 18
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 16.666666666666
$ws.Range("F18").Value = 32
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = 18.518518518518
$ws.Range("I18").Value = 81
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = 80
$ws.Range("L18").Value = 32.786885245901
$ws.Range("M18").Value = -4.705882352941
$ws.Range("N18").Value = -84.859813084112
# Row 19
$ws.Range("C19").Value = 31
$ws.Range("D19").Value = 30
$ws.Range("E19").Value = 3.333333333333
$ws.Range("F19").Value = 97
$ws.Range("G19").Value = 130
$ws.Range("H19").Value = -25.384615384615
$ws.Range("I19").Value = 195
$ws.Range("J19").Value = 272
$ws.Range("K19").Value = -28.308823529411
$ws.Range("L19").Value = -22
$ws.Range("M19").Value = 34.482758620689
$ws.Range("N19").Value = -22.310756972111
# Row 20
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = -85.714285714285
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 7.692307692307
$ws.Range("J20").Value = 37
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = -43.076923076923
$ws.Range("M20").Value = -26
$ws.Range("N20").Value = -95.595238095238
# Row 21
$ws.Range("C21").Value = 64
$ws.Range("D21").Value = 92
$ws.Range("E21").Value = -30.434782608695
$ws.Range("F21").Value = 235
$ws.Range("G21").Value = 298
$ws.Range("H21").Value = -21.140939597315
$ws.Range("I21").Value = 506
$ws.Range("J21").Value = 573
$ws.Range("K21").Value = -11.692844677137
$ws.Range("L21").Value = -14.237288135593
$ws.Range("M21").Value = 19.058823529411
$ws.Range("N21").Value = -75.590931017848
# Row 23
$ws.Range("F23").Value = 10
$ws.Range("H23").Value = 150
$ws.Range("I23").Value = 17
$ws.Range("J23").Value = 10
$ws.Range("K23").Value = 70
$ws.Range("L23").Value = -43.333333333333
$ws.Range("M23").Value = 183.333333333333
# Row 24
$ws.Range("C24").Value = 86
$ws.Range("D24").Value = 98
$ws.Range("E24").Value = -12.244897959183
$ws.Range("F24").Value = 342
$ws.Range("G24").Value = 371
$ws.Range("H24").Value = -7.816711590296
$ws.Range("I24").Value = 760
$ws.Range("J24").Value = 741
$ws.Range("K24").Value = 2.564102564102
$ws.Range("L24").Value = 15.501519756838
$ws.Range("M24").Value = 42.056074766355
# Row 25
$ws.Range("C25").Value = 48
$ws.Range("D25").Value = 66
$ws.Range("E25").Value = -27.272727272727
$ws.Range("F25").Value = 198
$ws.Range("G25").Value = 239
$ws.Range("H25").Value = -17.154811715481
$ws.Range("I25").Value = 459
$ws.Range("J25").Value = 409
$ws.Range("K25").Value = 12.224938875305
$ws.Range("L25").Value = 45.714285714285
# Row 26
$ws.Range("C26").Value = 36
$ws.Range("D26").Value = 48
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 116
$ws.Range("G26").Value = 151
$ws.Range("H26").Value = -23.178807947019
$ws.Range("I26").Value = 277
$ws.Range("J26").Value = 274
$ws.Range("K26").Value = 1.094890510948
$ws.Range("L26").Value = 1.838235294117
$ws.Range("M26").Value = 0.362318840579
# Row 27
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = -66.666666666666
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 13
$ws.Range("J27").Value = 20
$ws.Range("K27").Value = -35
$ws.Range("L27").Value = 44.444444444444
# Row 28
$ws.Range("C28").Value = 9
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = 125
$ws.Range("F28").Value = 24
$ws.Range("H28").Value = 60
$ws.Range("I28").Value = 41
$ws.Range("J28").Value = 29
$ws.Range("K28").Value = 41.379310344827
$ws.Range("L28").Value = 51.851851851851
# Row 29
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 4
$ws.Range("K29").Value = -75
$ws.Range("N29").Value = -93.75
# Row 30
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 4
$ws.Range("K30").Value = -75
$ws.Range("N30").Value = -92.857142857142
# Row 31
$ws.Range("G31").Value = 1
# Row 33
$ws.Range("L33").Value = -66.666666666666

# --- Cells changing type (number <-> text) need format fix-up ---

$ws.Range("C15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D29").Value = 1
$ws.Range("F29").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E29").Value = -100
$ws.Range("H29").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D30").Value = 1
$ws.Range("F30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E30").Value = -100
$ws.Range("H30").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

